# Adds 6 new OTP feedback rows (338-343) to the tracking sheet, continuing on
# from the existing last row (337), matching formatting of the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 337

$newRows = @(
    @{ A = 336; B = 42807; C = "transittracker@trimet.org"; D = "Complaint"; E = "Unhappy with trip plan; Confused by route map"; F = "Maximum walk distance too low"; H = "https://trimet.org/schedules/img/017.png"; I = "Weekday" },
    @{ A = 337; B = 42805; C = "triptech@trimet.org";       D = "Complaint"; E = "Thinks trip plan is incorrect"; F = "User accidentally selected PM instead of AM"; H = "http://trimet.org/#/planner/results/itin_num=3&from=Current"; I = "Saturday" },
    @{ A = 338; B = 42805; C = "triptech@trimet.org";       D = "Complaint"; E = "Dislikes walking directions"; H = "http://trimet.org/#/planner/results/itin_num=1&from=930%20NW%20NAITO%20PKWY,%20Portland::45.529015,-122.673416&to=Benson%20Polytechnic%20High,%20Portland::45.527313,-122.65263&Walk=1260&Arr=D"; I = "Saturday" },
    @{ A = 339; B = 42806; C = "triptech@trimet.org";       D = "Complaint"; E = "Unhappy with trip plan"; F = "Maximum walk distance too low"; H = "http://trimet.org/#planner/results/from=SW+ALLEN+BLVD+%26+SW+HALL+BLVD%3A%3A45.476464%2C-122.805481&to=11361+SW+LEVETON+DR%3A%3A45.386046%2C-122.794318&m=pm&walk=1260&arr=A"; I = "Weekday" },
    @{ A = 340; B = 42809; C = "triptech@trimet.org";       D = "Complaint"; E = "Dislikes inset maps"; H = "https://trimet.org/#planner/results/from=1511+SW+PARK+AVE%2C+Portland%3A%3A45.514206%2C-122.68472&to=631+NE+102ND+AVE%2C+Portland%3A%3A45.527378%2C-122.55855&m=pm&walk=1260&optimize=TRANSFERS&arr=A"; I = "Weekday" },
    @{ A = 341; B = 42811; C = "triptech@trimet.org";       D = "Complaint"; E = "Unhappy with trip plan"; F = "Maximum walk distance too low"; H = "https://trimet.org/#planner/results/from=187+SE+18TH+AVE%2C+Hillsboro%3A%3A45.521034%2C-122.96407&to=6360+SE+ALEXANDER+ST%2C+Hillsboro%3A%3A45.496017%2C-122.916664&m=pm&walk=1260&arr=A"; I = "Saturday" }
)

$r = $lastRow
foreach ($row in $newRows) {
    $r = $r + 1

    # Column A (record number) - reuse the bold/bordered/centered style used by
    # the rest of the "#" column.
    $ws.Cells.Item($lastRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row.A

    # Column B (date received) - reuse the date number-format style.
    $ws.Cells.Item($lastRow, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $row.B

    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($row.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    if ($row.ContainsKey("H")) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    if ($row.ContainsKey("I")) {
        $ws.Cells.Item($r, 9).Value = $row.I
    }
}
